$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 0.3333333333333333
$ws.Range("F2").Value = 0.7219999951205559
$ws.Range("G2").Value = 0.07200008696040985

# Row 3
$ws.Range("E3").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.2050000025153011
$ws.Range("G3").Value = 0.6490000112262161

# Row 4
$ws.Range("E4").Value = 0.3333333333333333
$ws.Range("F4").Value = 0.073000002364143
$ws.Range("G4").Value = 0.2789999018133741

# Row 5
$ws.Range("C5").Value = 0.4302324239064163
$ws.Range("D5").Value = 0.4382825810072314

# Row 6
$ws.Range("C6").Value = 0.3247054900642369
$ws.Range("D6").Value = 0.1234348381961154

# Row 7
$ws.Range("C7").Value = 0.2450620860293468
$ws.Range("D7").Value = 0.4382825807966532
